$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column cells to Text format so numeric-looking strings
# (e.g. "1.000", "0.9999") keep their exact text instead of being
# parsed into floating point numbers and losing trailing zeros.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Price column updates
$ws.Range("D2").Value = '27.089.50'
$ws.Range("D3").Value = '1.890.35'
$ws.Range("D4").Value = '0.9999'
$ws.Range("D5").Value = '307.52'
$ws.Range("D7").Value = '0.5140'
$ws.Range("D8").Value = '0.3743'
$ws.Range("D9").Value = '0.07213'
$ws.Range("D11").Value = '0.9060'
$ws.Range("D12").Value = '0.07643'
$ws.Range("D13").Value = '1.889.75'
$ws.Range("D14").Value = '94.83'
$ws.Range("D16").Value = '1.000'
$ws.Range("D17").Value = '0.000008486'
$ws.Range("D20").Value = '27.119.90'
$ws.Range("D21").Value = '5.075'
$ws.Range("D22").Value = '2.125.36'
$ws.Range("D24").Value = '6.411'
$ws.Range("D25").Value = '146.16'
$ws.Range("D26").Value = '1.790'
$ws.Range("D28").Value = '18.07'
$ws.Range("D29").Value = '114.62'
$ws.Range("D30").Value = '4.956'
$ws.Range("D31").Value = '4.854'
$ws.Range("D32").Value = '0.09190'
$ws.Range("D33").Value = '0.05082'
$ws.Range("D34").Value = '1.238'
$ws.Range("D36").Value = '2.973'
$ws.Range("D37").Value = '3.285'
$ws.Range("D38").Value = '2.603'
$ws.Range("D39").Value = '0.02001'
$ws.Range("D40").Value = '0.5607'
$ws.Range("D41").Value = '1.077'
$ws.Range("D42").Value = '6.667'
$ws.Range("D43").Value = '8.966'
$ws.Range("D44").Value = '118.06'
$ws.Range("D45").Value = '0.1516'
$ws.Range("D46").Value = '0.4807'
$ws.Range("D47").Value = '10.22'
$ws.Range("D48").Value = '1.0000'
$ws.Range("D50").Value = '37.47'
$ws.Range("D51").Value = '63.95'

# Volume(1h) / Coin / Link column updates
$ws.Range("E2").Value = '  +0.62%  '
$ws.Range("E3").Value = '  +1.48%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("E5").Value = '  +0.81%  '
$ws.Range("E7").Value = '  +1.39%  '
$ws.Range("E8").Value = '  +3.31%  '
$ws.Range("E9").Value = '  +0.56%  '
$ws.Range("E10").Value = '  +2.17%  '
$ws.Range("E11").Value = '  +1.07%  '
$ws.Range("E12").Value = '  +2.65%  '
$ws.Range("E13").Value = '  +1.86%  '
$ws.Range("E14").Value = '  +2.42%  '
$ws.Range("E15").Value = '  +0.65%  '
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("E17").Value = '  -0.04%  '
$ws.Range("E18").Value = '  +2.16%  '
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("E20").Value = '  +0.62%  '
$ws.Range("E21").Value = '  +0.92%  '
$ws.Range("E22").Value = '  +1.77%  '
$ws.Range("E23").Value = '  +2.21%  '
$ws.Range("E24").Value = '  -0.31%  '
$ws.Range("E25").Value = '  -1.33%  '
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("E28").Value = '  +1.19%  '
$ws.Range("E29").Value = '  +1.17%  '
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("E30").Value = '  +5.95%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("E31").Value = '  +3.86%  '
$ws.Range("E32").Value = '  -0.73%  '
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("E34").Value = '  +7.43%  '
$ws.Range("E35").Value = '  +2.63%  '
$ws.Range("E36").Value = '  -0.49%  '
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("E38").Value = '  +3.33%  '
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("E40").Value = '  +1.91%  '
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("E42").Value = '  +2.70%  '
$ws.Range("E43").Value = '  +4.86%  '
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("E45").Value = '  +3.11%  '
$ws.Range("E46").Value = '  +2.55%  '
$ws.Range("E47").Value = '  +1.39%  '
$ws.Range("E48").Value = '  +0.05%  '
$ws.Range("E49").Value = '  +1.74%  '
$ws.Range("E50").Value = '  +1.29%  '
$ws.Range("E51").Value = '  +1.47%  '
